$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all"
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# New daily row (2020/4/17 data)
$wsAll.Range("A9").Copy()
$wsAll.Range("A10").PasteSpecial(-4122)
$wsAll.Range("A10").Value = 43938
$wsAll.Range("B10").Value = 166
$wsAll.Range("C10").Value = 148
$wsAll.Range("D10").Value = 101
$wsAll.Range("E10").Value = 93
$wsAll.Range("F10").Value = 8
$wsAll.Range("G10").Value = 2
$wsAll.Range("H10").Value = 38

# Footnotes updated
$wsAll.Range("B11").Value = "※　24・34・53・58・59・60・158・161・163例目は市外在住者です。"
$wsAll.Range("B12").Value = "※　18件調査中"

# ---------------------------------------------------------------------
# Sheet "kobe"
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

# New daily row (2020/4/17 data)
$wsKobe.Range("A64").Copy()
$wsKobe.Range("A65").PasteSpecial(-4122)
$wsKobe.Range("A65").Value = 43938
$wsKobe.Range("B63").Copy()
$wsKobe.Range("B65").PasteSpecial(-4122)
$wsKobe.Range("B65").Value = 10
$wsKobe.Range("C65").Value = 1097
$wsKobe.Range("D65").Value = 5
$wsKobe.Range("E65").Value = 166
$wsKobe.Range("F65").Value = 101
$wsKobe.Range("G65").Value = 93
$wsKobe.Range("H65").Value = 8
$wsKobe.Range("I65").Value = 2
$wsKobe.Range("J65").Value = 36

# New footnote about out-of-city residents
$wsKobe.Range("E66").Value = "※　24・34・53・58・59・60・158・161・163例目は市外在住者です。"

# ---------------------------------------------------------------------
# Sheet "other"
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

# Fill in the previously-blank daily row (2020/4/17 data)
$wsOther.Range("A40").Value = 43938
$wsOther.Range("B40").Value = 0
$wsOther.Range("C40").Value = 9
$wsOther.Range("D40").Value = 7
$wsOther.Range("E40").Value = 6
$wsOther.Range("F40").Value = 1
$wsOther.Range("G40").Value = 0
$wsOther.Range("H40").Value = 2

# ---------------------------------------------------------------------
# Selections / scroll positions (match the updated cursor state)
# ---------------------------------------------------------------------
$wsKobe.Range("G70").Select()
$wsOther.Range("D41").Select()
$wsAll.Range("B12").Select()
